# Update database and change read_price algorithm
# Rolling 12-month-period table: drop the oldest period (1396/12) and
# shift every period one column to the left, adding the newest period
# (1401/12) in the rightmost column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two period-header rows (row 8 and row 24), columns E..I ---
$headers = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)

for ($i = 0; $i -lt 5; $i++) {
    $col = 5 + $i  # E=5 .. I=9
    $ws.Cells.Item(8, $col).Value = $headers[$i]
    $ws.Cells.Item(24, $col).Value = $headers[$i]
}

# --- Shift the numeric data one column to the left (drop oldest year,
#     append newly computed value for the newest year) for each data row ---

function Shift-Row($row, $newLast) {
    $e = $ws.Cells.Item($row, 5).Value2
    $f = $ws.Cells.Item($row, 6).Value2
    $g = $ws.Cells.Item($row, 7).Value2
    $h = $ws.Cells.Item($row, 8).Value2
    $i = $ws.Cells.Item($row, 9).Value2

    $ws.Cells.Item($row, 5).Value = $f
    $ws.Cells.Item($row, 6).Value = $g
    $ws.Cells.Item($row, 7).Value = $h
    $ws.Cells.Item($row, 8).Value = $i
    $ws.Cells.Item($row, 9).Value = $newLast
}

Shift-Row 10 10729
Shift-Row 11 0
Shift-Row 12 0
Shift-Row 13 9673
Shift-Row 14 9875
Shift-Row 15 23157
Shift-Row 16 23774
Shift-Row 17 399681
Shift-Row 18 0
Shift-Row 19 130495
Shift-Row 20 607384

Shift-Row 26 152
Shift-Row 27 829
